$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados" timestamp string in cell A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Septiembre de 2020 a las 09:58"

# Row 55 - Singapur
$ws.Range("B55").Value = 57454
$ws.Range("C55").Value = 48
$ws.Range("E55").Value = 663

# Row 62 - Armenia
$ws.Range("B62").Value = 45969
$ws.Range("C62").Value = 107
$ws.Range("D62").Value = 41693
$ws.Range("E62").Value = 3357
$ws.Range("G62").Value = 3
$ws.Range("H62").Value = 919

# Row 92 - Hungria
$ws.Range("B92").Value = 13153
$ws.Range("C92").Value = 844
$ws.Range("D92").Value = 4117
$ws.Range("E92").Value = 8394
$ws.Range("G92").Value = 5
$ws.Range("H92").Value = 642

# Row 143 - Estonia
$ws.Range("B143").Value = 2698
$ws.Range("C143").Value = 22
$ws.Range("E143").Value = 378

# Row 160 - Letonia
$ws.Range("B160").Value = 1477
$ws.Range("C160").Value = 3
$ws.Range("E160").Value = 194
